$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting of the existing header cells (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the new data columns I (I0) and J (IF) for rows 2-9
$data = @{
    2 = @(5, 5)
    3 = @(7, 7)
    4 = @(8, 9)
    5 = @(3, 3)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(8, 8)
    9 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
